# Croatia 3NL - "Atualização de bases das ligas" update.
#
# The upstream odds feed re-synced some fixtures that share the exact same
# kick-off date, causing their rows (everything except the running index in
# column A) to trade places in the sheet. This script reproduces that by
# rotating the contents of columns B:AC among the affected rows, leaving
# column A (the original sequential id) untouched on every row.
#
# Each inner array below is a cycle of 1-based worksheet row numbers whose
# B:AC payloads must be rotated one position "to the left", i.e. row
# cycle[0] receives what used to be in cycle[1], cycle[1] receives what
# used to be in cycle[2], ..., and the last one wraps back to cycle[0].
# Most cycles are simple 2-row swaps, but a couple of groups of rows moved
# around together (a 4-row and a 3-row rotation).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowCycles = @(
  @(8, 9),
  @(14, 15),
  @(18, 19, 20, 21),
  @(22, 24),
  @(29, 30),
  @(33, 34),
  @(35, 36),
  @(41, 42),
  @(43, 44),
  @(50, 51),
  @(58, 59),
  @(65, 66, 67),
  @(73, 74),
  @(76, 77),
  @(83, 84),
  @(93, 94),
  @(101, 102),
  @(107, 108)
)

foreach ($cycle in $rowCycles) {
    $n = $cycle.Length

    # Snapshot the B:AC payload of every row in the cycle before writing
    # anything back, so rows that feed into each other don't clobber data
    # that hasn't been read yet.
    $snapshots = @()
    for ($i = 0; $i -lt $n; $i++) {
        $r = $cycle[$i]
        $snapshots += ,($ws.Range("B$r`:AC$r").Value())
    }

    for ($i = 0; $i -lt $n; $i++) {
        $r = $cycle[$i]
        $sourceIndex = ($i + 1) % $n
        $ws.Range("B$r`:AC$r").Value = $snapshots[$sourceIndex]
    }
}
